# Update the IGCC Netting Flows model for MM&MV:
# - shift the timestamp series forward by one day (2026-02-16 -> 2026-02-17)
# - refresh the Import/Export MW readings for each 15-minute interval
# - the new series has one fewer row (30 data rows instead of 31), so the
#   trailing row is removed afterwards

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46070,                 1.868,  0),
    @(46070.01041666666,     3.637,  0),
    @(46070.02083333334,     1.761,  0),
    @(46070.03125,           2.187,  0.095),
    @(46070.04166666666,     1.902,  0),
    @(46070.05208333334,     8.786,  0),
    @(46070.0625,            1.81,   0),
    @(46070.07291666666,     2.175,  0.9370000000000001),
    @(46070.08333333334,     12.251, 0),
    @(46070.09375,           24.604, 0.06),
    @(46070.10416666666,     36.973, 0),
    @(46070.11458333334,     32.228, 0),
    @(46070.125,             14.386, 0),
    @(46070.13541666666,     20.857, 0),
    @(46070.14583333334,     21.028, 0),
    @(46070.15625,           16.087, 0),
    @(46070.16666666666,     9.474,  0),
    @(46070.17708333334,     38.836, 0),
    @(46070.1875,            26.323, 0),
    @(46070.19791666666,     0.648,  1.619),
    @(46070.20833333334,     0,      20.901),
    @(46070.21875,           0,      12.428),
    @(46070.22916666666,     5.941,  0.407),
    @(46070.23958333334,     6.279,  0.383),
    @(46070.25,              0.077,  3.59),
    @(46070.26041666666,     14.405, 0),
    @(46070.27083333334,     31.983, 0),
    @(46070.28125,           37.544, 0),
    @(46070.29166666666,     22.599, 0),
    @(46070.3125,            0,      0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# The refreshed series only has 30 rows of data (was 31), drop the old last row
$ws.Rows.Item(32).Delete()
